$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new data row (row 51) with the same layout as preceding rows.
$row = 51
$ws.Cells.Item($row, 1).Value = 0
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0.363636
$ws.Cells.Item($row, 6).Value = -2.446380104769124
$ws.Cells.Item($row, 7).Value = "query"
